# Applies a re-ordering of the species-observation rows 16-21 on the
# "Artfynd" sheet. The columns that vary between rows (A, B, D, E, F, G,
# H, Q, R, Z, AB) are re-shuffled according to the mapping observed in
# the target diff; all other columns are identical across these rows and
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that differ row-to-row and need to be moved around.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "Z", "AB")

$rows = @(16, 17, 18, 19, 20, 21)

# Snapshot the current ("before") values for every relevant cell so that
# writes below never read back an already-modified value.
$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowData
}

# Maps each destination row to the source row whose data it should now
# contain (derived from the committed change).
$rowMap = @{
    16 = 21
    17 = 19
    18 = 20
    19 = 18
    20 = 17
    21 = 16
}

foreach ($destRow in $rows) {
    $srcRow = $rowMap[$destRow]
    $srcData = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $srcData[$col]
    }
}
